$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# fix(attendance-import): fix import attendance
# Rename the "employee_code" header to "nik"
$ws.Range("C1").Value = "nik"

# Append a sample attendance record in row 2.
# Use a leading apostrophe on the numeric/date/time-looking values so Excel
# keeps them as text (matching the source import's string-typed columns)
# instead of silently coercing them into numbers / date serials.
$ws.Range("A2").Value = "WCK760"
$ws.Range("B2").Value = "SIF141"
$ws.Range("C2").Value = "'123123"
$ws.Range("D2").Value = "'2023-10-01"
$ws.Range("E2").Value = "'08:00"
$ws.Range("F2").Value = "'18:00"
$ws.Range("G2").Value = "on time"
